# Cambio_Pantalla_Final_HyO_insumo_ConPagos.xlsx
# "finalizada para cambios de retirar recibo"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # dataConPagos
$ws2 = $wb.Worksheets.Item(2)   # aa

# --- 1. Flip the "CAMBIO" (SI -> NO) flag for most rows, keeping a handful as SI ---
$rowRanges = @(
    @(2, 8),
    @(10, 11),
    @(13, 23),
    @(29, 38),
    @(40, 52),
    @(55, 80),
    @(90, 101),
    @(120, 123)
)

foreach ($range in $rowRanges) {
    $start = $range[0]
    $end = $range[1]
    for ($r = $start; $r -le $end; $r++) {
        $ws1.Cells.Item($r, 3).Value = "NO"
    }
}

# --- 2. For the last block (rows 120-123) the "ID CAMBIO" column switches from the
#        numeric 385 to the text value "65" (stored as text, hence the later
#        numberStoredAsText ignored-error on D120:D123) ---
for ($r = 120; $r -le 123; $r++) {
    $ws1.Range("D$r").NumberFormat = "@"
    $ws1.Range("D$r").Value = "65"
}

# --- 3. Sheet/window state: "dataConPagos" becomes the active tab/sheet, scrolled
#        down with C101 selected; "aa" loses its tabSelected flag (simply by no
#        longer being the active sheet) ---
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 91
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("C101").Select()
